$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Simulacion_1")
$ws.Range("B2").Value = 85.33359125119999
$ws.Range("C2").Value = 16.121362785
$ws.Range("D2").Value = 5.293199612789435
$ws.Range("E2").Value = 12.70367907069464
$ws.Range("I2").Value = 4.20404357200725
$ws.Range("K2").Value = 6.13370815875
$ws.Range("L2").Value = 16.35880434125
$ws.Range("M2").Value = 1.92966458674275
$ws.Range("Q2").Value = 1.761901158124151
$ws.Range("S2").Value = 2.801222578933642
$ws.Range("T2").Value = 0.0747009981676096
$ws.Range("U2").Value = 5.065529323443576

$ws = $wb.Worksheets.Item("Simulacion_2")
$ws.Range("B2").Value = 95.78052290468307
$ws.Range("C2").Value = 19.82325696100919
$ws.Range("D2").Value = 4.831724831750702
$ws.Range("E2").Value = 11.59613959620168
$ws.Range("I2").Value = 4.878013028644973
$ws.Range("J2").Value = 26.09839007200725
$ws.Range("K2").Value = 7.117030972636377
$ws.Range("L2").Value = 18.98135909937088
$ws.Range("M2").Value = 2.239017943991404
$ws.Range("Q2").Value = 1.740288072173941
$ws.Range("R2").Value = 0.8081863586270529
$ws.Range("S2").Value = 2.766860228875015
$ws.Range("T2").Value = 0.07378464761836426
$ws.Range("U2").Value = 5.003390922463434

$ws = $wb.Worksheets.Item("Simulacion_3")
$ws.Range("B2").Value = 87.40324667284139
$ws.Range("C2").Value = 15.77147756935734
$ws.Range("D2").Value = 5.541855307371997
$ws.Range("E2").Value = 13.30045273769279
$ws.Range("G2").Value = 1.4709
$ws.Range("H2").Value = 23.0441
$ws.Range("I2").Value = 5.218882496783531
$ws.Range("J2").Value = 27.92211302864497
$ws.Range("K2").Value = 7.614360222911484
$ws.Range("L2").Value = 20.30775280573349
$ws.Range("M2").Value = 2.395477726127953
$ws.Range("O2").Value = 22.716
$ws.Range("P2").Value = 1.103234042553192
$ws.Range("Q2").Value = 2.615273423367069
$ws.Range("R2").Value = 1.21452783515595
$ws.Range("S2").Value = 4.157987484054254
$ws.Range("T2").Value = 0.1108822332659908
$ws.Range("U2").Value = 7.519005339092387

